$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $pPrXml, $innerXml) {
    $p = $d.Paragraphs($paraIndex).Range
    $xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrXml + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $p.InsertXML($xmlFragment)
}

Replace-ParagraphXml 4 '<w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' '<w:r><w:tab/><w:t>To display information about</w:t></w:r><w:r><w:t xml:space="preserve"> all</w:t></w:r><w:r><w:t xml:space="preserve"> module</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> found in a</w:t></w:r><w:r><w:t xml:space="preserve"> directory using its</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t xml:space="preserve">xml </w:t></w:r><w:r><w:t>files and be able to save an enable/disable for each module</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r>'
Replace-ParagraphXml 8 '<w:pPr><w:ind w:left="1440"/></w:pPr>' '<w:r><w:t xml:space="preserve">All public methods of a module will be displayed under the module name in </w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TreeView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> control</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>The user can select each method</w:t></w:r><w:r><w:t xml:space="preserve"> which will then display all information about that module and the selected method</w:t></w:r><w:r><w:t xml:space="preserve"> in a different field</w:t></w:r><w:r><w:t>.</w:t></w:r>'
Replace-ParagraphXml 10 '<w:pPr><w:ind w:left="1440"/></w:pPr>' '<w:r><w:t>Each module</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">will have </w:t></w:r><w:r><w:t>its</w:t></w:r><w:r><w:t xml:space="preserve"> own enable/disable check box to indicate whether the module </w:t></w:r><w:r><w:t>will be</w:t></w:r><w:r><w:t xml:space="preserve"> used or not. </w:t></w:r><w:r><w:t xml:space="preserve">It will have a “load </w:t></w:r><w:r><w:t>modul</w:t></w:r><w:r><w:t>es</w:t></w:r><w:r><w:t xml:space="preserve">” button to allow the user to select the location of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dll’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
Replace-ParagraphXml 12 '<w:pPr><w:ind w:left="1440"/></w:pPr>' '<w:r><w:t>There will be a “save configuration” button which will save the state of each modules enable/disable check box</w:t></w:r><w:r><w:t xml:space="preserve">. On </w:t></w:r><w:r><w:t xml:space="preserve">application </w:t></w:r><w:r><w:t xml:space="preserve">startup, </w:t></w:r><w:r><w:t>if the config file exists</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t xml:space="preserve">the previously saved settings </w:t></w:r><w:r><w:t>will be</w:t></w:r><w:r><w:t xml:space="preserve"> loaded. The last selected directory will also be saved so all modules in that directory will be loaded and displayed</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r>'
Replace-ParagraphXml 15 '<w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Breakdown</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r>'
Replace-ParagraphXml 16 '<w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' '<w:r><w:tab/><w:t xml:space="preserve">1.)  A </w:t></w:r><w:r><w:t>“load modules”</w:t></w:r><w:r><w:t xml:space="preserve"> button will allow the user to select a directory where </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> files are </w:t></w:r><w:r><w:t xml:space="preserve">located. The program will check that the directory contains at least one </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file.</w:t></w:r>'
Replace-ParagraphXml 18 '<w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' '<w:r><w:tab/><w:t xml:space="preserve">2.)  The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TreeView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> area will display each module (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dll</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file) and its methods in child nodes. Each module will have an associated enable/disable check box.</w:t></w:r>'
Replace-ParagraphXml 20 '<w:pPr><w:ind w:left="1440" w:hanging="1440"/></w:pPr>' '<w:r><w:tab/><w:t xml:space="preserve">3.)  The user can select any of the methods from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TreeView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> area. The information about the selected method will be displayed in a </w:t></w:r><w:r><w:t>separate field</w:t></w:r><w:r><w:t xml:space="preserve"> next to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TreeView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> area.</w:t></w:r>'

Write-Host "Done."
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "$i : " $d.Paragraphs($i).Range.Text
}
